$d = $word.ActiveDocument

function Replace-RangeWithXml($oldText, $pkgXml) {
    $full = $d.Content.Text
    $idx = $full.IndexOf($oldText)
    if ($idx -lt 0) {
        throw "Target text not found: $($oldText.Substring(0, [Math]::Min(40, $oldText.Length)))"
    }
    $rng = $d.Range($idx, $idx + $oldText.Length)
    $rng.InsertXML($pkgXml) | Out-Null
}

$oldText1 = @'
Objetivo GeralPermitir aos estudantes que compreendam os mecanismos de obtenção da influencia de diversos fatores (variáveis independentes de um processo) sobre as variáveis resposta (dependentes), através da análise multivariada.Objetivos EspecíficosSaber planejar e executar um experimento fatorial completo e fracionadoSaber analisar os resultados propondo a condição de melhor ajuste que otimiza os valores da variável resposta na região experimental estudadaDominar, pelo menos, um software comercial sobre o assuntoSaber modelar um processo, com base em dados empíricos
'@

$pkgXml1 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Objetivo Geral</w:t><w:br/><w:t>Permitir aos estudantes que compreendam os mecanismos de obtenção da influencia de diversos fatores (variáveis independentes de um processo) sobre as variáveis resposta (dependentes), através da análise multivariada.</w:t><w:br/><w:br/><w:t>Objetivos Específicos</w:t><w:br/><w:t>Saber planejar e executar um experimento fatorial completo e fracionado</w:t><w:br/><w:t>Saber analisar os resultados propondo a condição de melhor ajuste que otimiza os valores da variável resposta na região experimental estudada</w:t><w:br/><w:t>Dominar, pelo menos, um software comercial sobre o assunto</w:t><w:br/><w:t>Saber modelar um processo, com base em dados empíricos</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$oldText2 = @'
1. MONTGOMERY, D.C., Design and Analysis of Experiments, Wiley, 19912. BOX, G.E.; HUNTER, W.G.; HUNTER, J.S., Statistic for Experimenters, John Wiley & Sons, New York, 1978. 3. TAGUCHI, G.; WU, YU-IN., Introduction to off-Line Quality Control. Central Japan Quality Control Association. Meieki Nakamura-Ku Magaya, Japan, 1979. 4. BRUNS, R.E., Como Fazer Experimentos, Editora UNICAMP, 2010. 5. COX, D.R., Planning of Experiments, Wiley 1976. 6. COX, G.M.; COCHRAN, W.G., Experimental Desing. Wiley 1976. 7. SILVA M.B. et al, Design of Experiments-Applications, Editora Intech, 2013
'@

$pkgXml2 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>1. MONTGOMERY, D.C., Design and Analysis of Experiments, Wiley, 1991</w:t><w:br/><w:t xml:space="preserve">2. BOX, G.E.; HUNTER, W.G.; HUNTER, J.S., Statistic for Experimenters, John Wiley &amp; Sons, New York, 1978. </w:t><w:br/><w:t xml:space="preserve">3. TAGUCHI, G.; WU, YU-IN., Introduction to off-Line Quality Control. Central Japan Quality Control Association. Meieki Nakamura-Ku Magaya, Japan, 1979. </w:t><w:br/><w:t xml:space="preserve">4. BRUNS, R.E., Como Fazer Experimentos, Editora UNICAMP, 2010. </w:t><w:br/><w:t xml:space="preserve">5. COX, D.R., Planning of Experiments, Wiley 1976. </w:t><w:br/><w:t xml:space="preserve">6. COX, G.M.; COCHRAN, W.G., Experimental Desing. Wiley 1976. </w:t><w:br/><w:t>7. SILVA M.B. et al, Design of Experiments-Applications, Editora Intech, 2013</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

Replace-RangeWithXml $oldText1 $pkgXml1
Replace-RangeWithXml $oldText2 $pkgXml2

Write-Output "Edit complete"
